$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the header row: a new leading " Invoice Code" column replaces the
# old trailing "Amount w/o Tax" / "Tax Value" / "Total included tax" /
# "Certified Amount" columns with just two renamed "Certified ..." columns,
# and the now-unused last column is removed entirely.
$ws.Range("A1").Value = " Invoice Code"
$ws.Range("B1").Value = " Invoice number"
$ws.Range("C1").Value = "(*)Vendor"
$ws.Range("D1").Value = "(*)Budget Head"
$ws.Range("E1").Value = "(*)Invoice date"
$ws.Range("F1").Value = "(*)Project"
$ws.Range("G1").Value = "(*)Description of Services"
$ws.Range("H1").Value = "Certified Amount w/o Tax"
$ws.Range("I1").Value = "Certified Tax Amount"

# Drop the old trailing column (used to hold "Certified Amount").
$ws.Columns("J:J").Delete()

# Match the author's new selection.
$ws.Range("G3").Select()
